$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric (e.g. "1.00", "5.57") must be
# forced to Text format first, otherwise Excel auto-converts them to numbers
# (matching how the source data is stored as literal text in the workbook).
# NumberFormat is set per-cell (not as one multi-area range) since multi-area
# Range.NumberFormat assignment only reliably applies to the first area.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

$ws.Range('D2').Value = '65.338.63'
$ws.Range('E2').Value = '  -0.86%  '
$ws.Range('D3').Value = '3.331.06'
$ws.Range('E3').Value = '  -4.46%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '572.78'
$ws.Range('E5').Value = '  -1.72%  '
$ws.Range('D6').Value = '177.06'
$ws.Range('E6').Value = '  +2.68%  '
$ws.Range('D7').Value = '0.615'
$ws.Range('E7').Value = '  +2.97%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '3.329.83'
$ws.Range('E9').Value = '  -4.32%  '
$ws.Range('E10').Value = '  -2.24%  '
$ws.Range('E11').Value = '  -0.29%  '
$ws.Range('D12').Value = '0.406'
$ws.Range('E12').Value = '  -1.11%  '
$ws.Range('D13').Value = '3.907.36'
$ws.Range('E13').Value = '  -4.27%  '
$ws.Range('E14').Value = '  +0.06%  '
$ws.Range('E15').Value = '  -5.16%  '
$ws.Range('D16').Value = '65.371.81'
$ws.Range('E16').Value = '  -0.91%  '
$ws.Range('E17').Value = '  -1.86%  '
$ws.Range('D18').Value = '3.339.88'
$ws.Range('E18').Value = '  -3.92%  '
$ws.Range('D19').Value = '5.73'
$ws.Range('E19').Value = '  -3.14%  '
$ws.Range('D20').Value = '13.33'
$ws.Range('E20').Value = '  -4.12%  '
$ws.Range('D21').Value = '361.01'
$ws.Range('E21').Value = '  -1.59%  '
$ws.Range('D22').Value = '7.41'
$ws.Range('E22').Value = '  -4.25%  '
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').Value = '70.98'
$ws.Range('E24').Value = '  -2.43%  '
$ws.Range('D25').Value = '0.515'
$ws.Range('E25').Value = '  -3.55%  '
$ws.Range('E26').Value = '  -4.32%  '
$ws.Range('D27').Value = '9.50'
$ws.Range('E27').Value = '  -1.04%  '
$ws.Range('E28').Value = '  -1.07%  '
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('E30').Value = '  -1.81%  '
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  -0.05%  '
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').Value = '5.57'
$ws.Range('E32').Value = '  -3.28%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').Value = '22.84'
$ws.Range('E33').Value = '  -5.35%  '
$ws.Range('E34').Value = '  -4.75%  '
$ws.Range('E35').Value = '  -6.82%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '1.48'
$ws.Range('E36').Value = '  -3.44%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').Value = '159.93'
$ws.Range('E37').Value = '  -0.37%  '
$ws.Range('D38').Value = '0.844'
$ws.Range('E38').Value = '  -5.12%  '
$ws.Range('D39').Value = '27.43'
$ws.Range('E39').Value = '  -6.90%  '
$ws.Range('E40').Value = '  -0.96%  '
$ws.Range('D41').Value = '2.700.43'
$ws.Range('E41').Value = '  -4.47%  '
$ws.Range('D42').Value = '2.48'
$ws.Range('E42').Value = '  -3.12%  '
$ws.Range('E43').Value = '  -4.44%  '
$ws.Range('D44').Value = '4.26'
$ws.Range('E44').Value = '  -4.42%  '
$ws.Range('D45').Value = '39.74'
$ws.Range('E45').Value = '  -0.82%  '
$ws.Range('D46').Value = '0.0663'
$ws.Range('E46').Value = '  -2.76%  '
$ws.Range('D47').Value = '332.54'
$ws.Range('E47').Value = '  +2.34%  '
$ws.Range('D48').Value = '23.81'
$ws.Range('E48').Value = '  -1.16%  '
$ws.Range('D49').Value = '0.0277'
$ws.Range('E49').Value = '  -3.53%  '
$ws.Range('E50').Value = '  +1.65%  '
$ws.Range('E51').Value = '  +0.09%  '
